$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "modif 11H30"
$ws.Range("B2").Select()
